$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the logged time for E33 (added two more 21-minute entries) which
# cascades through the dependent SUM/ratio formulas automatically.
$ws.Range("E33").Formula = "=(1/60)*(9+13+21+21+21+21)"

# Move the active selection to reflect where the user left off editing.
$ws.Range("H33").Select()
